$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04202566666666666
$ws.Range("H2").Value = 0.126077
$ws.Range("I2").Value = 0.001003231639737821
$ws.Range("J2").Value = 0.001003231639737821
$ws.Range("M2").Value = 2.675565666666666
$ws.Range("N2").Value = 8.026696999999999
$ws.Range("O2").Value = 0.2572287714720184
$ws.Range("P2").Value = 0.2572287714720184
$ws.Range("Q2").Value = 0.1124424308521111
$ws.Range("R2").Value = 1.011981877669
$ws.Range("S2").Value = 0.0002580600421916182
$ws.Range("T2").Value = 0.0002580600421916182

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04202566666666666
$ws.Range("H3").Value = 0.126077
$ws.Range("I3").Value = 0.001003231639737821
$ws.Range("J3").Value = 0.001003231639737821
$ws.Range("O3").Value = 0.2087228937794146
$ws.Range("P3").Value = 0.2087228937794146
$ws.Range("Q3").Value = 0.09123905314611111
$ws.Range("R3").Value = 0.821151478315
$ws.Range("S3").Value = 0.0002093974109771451
$ws.Range("T3").Value = 0.0002093974109771452

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04202566666666666
$ws.Range("H4").Value = 0.126077
$ws.Range("I4").Value = 0.001003231639737821
$ws.Range("J4").Value = 0.001003231639737821
$ws.Range("M4").Value = 5.554905000000001
$ws.Range("N4").Value = 16.664715
$ws.Range("O4").Value = 0.5340483347485671
$ws.Range("P4").Value = 0.5340483347485669
$ws.Range("Q4").Value = 0.233448585895
$ws.Range("R4").Value = 2.101037273055
$ws.Range("S4").Value = 0.0005357741865690575
$ws.Range("T4").Value = 0.0005357741865690575

# Row 5
$ws.Range("I5").Value = 0.9479341900351343
$ws.Range("J5").Value = 0.9479341900351345
$ws.Range("M5").Value = 2.675565666666666
$ws.Range("N5").Value = 8.026696999999999
$ws.Range("O5").Value = 0.2572287714720184
$ws.Range("P5").Value = 0.2572287714720184
$ws.Range("Q5").Value = 106.2446800852819
$ws.Range("R5").Value = 956.2021207675368
$ws.Range("S5").Value = 0.2438359471390604
$ws.Range("T5").Value = 0.2438359471390605

# Row 6
$ws.Range("I6").Value = 0.9479341900351343
$ws.Range("J6").Value = 0.9479341900351345
$ws.Range("O6").Value = 0.2087228937794146
$ws.Range("P6").Value = 0.2087228937794146
$ws.Range("R6").Value = 775.8901640064951
$ws.Range("S6").Value = 0.1978555672565788
$ws.Range("T6").Value = 0.1978555672565788

# Row 7
$ws.Range("I7").Value = 0.9479341900351343
$ws.Range("J7").Value = 0.9479341900351345
$ws.Range("M7").Value = 5.554905000000001
$ws.Range("N7").Value = 16.664715
$ws.Range("O7").Value = 0.5340483347485671
$ws.Range("P7").Value = 0.5340483347485669
$ws.Range("Q7").Value = 220.581057673835
$ws.Range("R7").Value = 1985.229519064515
$ws.Range("S7").Value = 0.5062426756394952
$ws.Range("T7").Value = 0.5062426756394952

# Row 8
$ws.Range("G8").Value = 2.139026333333333
$ws.Range("H8").Value = 6.417078999999999
$ws.Range("I8").Value = 0.05106257832512778
$ws.Range("J8").Value = 0.05106257832512778
$ws.Range("M8").Value = 2.675565666666666
$ws.Range("N8").Value = 8.026696999999999
$ws.Range("O8").Value = 0.2572287714720184
$ws.Range("P8").Value = 0.2572287714720184
$ws.Range("Q8").Value = 5.723105417562554
$ws.Range("R8").Value = 51.50794875806299
$ws.Range("S8").Value = 0.01313476429076633
$ws.Range("T8").Value = 0.01313476429076633

# Row 9
$ws.Range("G9").Value = 2.139026333333333
$ws.Range("H9").Value = 6.417078999999999
$ws.Range("I9").Value = 0.05106257832512778
$ws.Range("J9").Value = 0.05106257832512778
$ws.Range("O9").Value = 0.2087228937794146
$ws.Range("P9").Value = 0.2087228937794146
$ws.Range("Q9").Value = 4.643893905500556
$ws.Range("R9").Value = 41.795045149505
$ws.Range("S9").Value = 0.01065792911185869
$ws.Range("T9").Value = 0.01065792911185869

# Row 10
$ws.Range("G10").Value = 2.139026333333333
$ws.Range("H10").Value = 6.417078999999999
$ws.Range("I10").Value = 0.05106257832512778
$ws.Range("J10").Value = 0.05106257832512778
$ws.Range("M10").Value = 5.554905000000001
$ws.Range("N10").Value = 16.664715
$ws.Range("O10").Value = 0.5340483347485671
$ws.Range("P10").Value = 0.5340483347485669
$ws.Range("Q10").Value = 11.882088074165
$ws.Range("R10").Value = 106.938792667485
$ws.Range("S10").Value = 0.02726988492250277
$ws.Range("T10").Value = 0.02726988492250276

